$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Filtered Stats")

# Updated counts after removing indels from the invariants (Num Invariants) column set
# Column B = Num SNPs, C = Num Multiallelic, D = Num Invariants (per diff order)
$ws.Range("B2").Value = 104205
$ws.Range("C2").Value = 7321
$ws.Range("D2").Value = 683440

$ws.Range("B3").Value = 103032
$ws.Range("C3").Value = 6933
$ws.Range("D3").Value = 736948

$ws.Range("B4").Value = 109859
$ws.Range("C4").Value = 8366
$ws.Range("D4").Value = 626870

$ws.Range("B5").Value = 82132
$ws.Range("C5").Value = 5771
$ws.Range("D5").Value = 538119

$ws.Range("B6").Value = 97517
$ws.Range("C6").Value = 7564
$ws.Range("D6").Value = 523321

$ws.Range("B7").Value = 123970
$ws.Range("C7").Value = 8191
$ws.Range("D7").Value = 955746

$ws.Range("B8").Value = 92452
$ws.Range("C8").Value = 6822
$ws.Range("D8").Value = 569203

$ws.Range("B9").Value = 146604
$ws.Range("C9").Value = 10603
$ws.Range("D9").Value = 891241

$ws.Range("B10").Value = 91501
$ws.Range("C10").Value = 6347
$ws.Range("D10").Value = 597460

$ws.Range("B11").Value = 96478
$ws.Range("C11").Value = 6575
$ws.Range("D11").Value = 683749

$ws.Range("B12").Value = 93288
$ws.Range("C12").Value = 6132
$ws.Range("D12").Value = 681067

$ws.Range("B13").Value = 4181
$ws.Range("C13").Value = 434
$ws.Range("D13").Value = 19068

# Update the active selection to reflect the last cell the author was working on
$ws.Activate()
$ws.Range("J12").Select()

$wb.Save()
